# Add tests for unit field batch overlay.
#
# The "TEST" sheet gains two new columns (Subject:topic FAST / Subject:genre
# FAST) inserted right after the existing "Subject:topic" column, and a new
# "unit"/"Unit" column appended at the end of the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns after column O ("Subject:topic"), i.e. at P:Q.
# Excel inherits formatting from the column being pushed right (old P),
# which is exactly the style the header/data rows already use at O.
$ws.Columns("P:Q").Insert()

$ws.Range("P1").Value = "Subject:topic FAST"
$ws.Range("Q1").Value = "Subject:genre FAST"
$ws.Range("P2").Value = "Subject:topic FAST"
$ws.Range("Q2").Value = "Subject:genre FAST"

# --- Append a new "unit" / "Unit" column at the end (now column AL), by
# copying the last existing column (AK) one step right and overwriting its
# values - this keeps the header/data formatting consistent with the rest
# of the row while giving the new column its own text.
$ws.Columns("AK:AK").Copy()
$ws.Columns("AL:AL").Insert()

$ws.Range("AL1").Value = "unit"
$ws.Range("AL2").Value = "Unit"

# --- Update the view: scroll right so the new columns are visible and
# select the newly added AL2 cell.
$ws.Range("AL2").Select()
$excel.ActiveWindow.ScrollColumn = 31
